$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for "^_^^_^" (currently row 34) needs to move down so it sits
# right above the "WvCly" row (originally row 45), i.e. it becomes the
# last entry of the "Ok"/"Verificar" block before WvCly.

# 1) Capture the values of row 34 (A..H) before touching anything.
$capturedRow = @()
for ($col = 1; $col -le 8; $col++) {
    $capturedRow += ,$ws.Cells.Item(34, $col).Value()
}

# 2) Delete row 34 entirely - rows 35-45 shift up to become rows 34-44.
$ws.Rows.Item(34).Delete()

# 3) After the delete, "WvCly" (originally row 45) now sits at row 44.
#    Insert a fresh blank row above it so it shifts back down to row 45,
#    leaving a new blank row 44 for the moved data.
$ws.Rows.Item(44).Insert()

# 4) Fill the new row 44 with the captured data, restoring it just above
#    WvCly.
for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item(44, $col).Value = $capturedRow[$col - 1]
}
